$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 12:50"

# Country rows whose name/order and/or case counts changed
$rowsData = @(
    @{ Row = 4; Vals = @("Estados Unidos", 188639, 109, 7251, 177329, 4576, 6, 4059) },
    @{ Row = 23; Vals = @("Noruega", 4699, 58, 13, 4646, 105, 1, 40) },
    @{ Row = 36; Vals = @("Japon", 2178, 0, 472, 1649, 69, 0, 57) },
    @{ Row = 42; Vals = @("Finlandia", 1446, 28, 10, 1419, 56, 0, 17) },
    @{ Row = 79; Vals = @("Azerbaiyan", 359, 61, 26, 328, 7, 0, 5) },
    @{ Row = 80; Vals = @("Moldavia", 353, 0, 22, 327, 44, 0, 4) },
    @{ Row = 81; Vals = @("Costa Rica", 347, 0, 4, 341, 8, 0, 2) },
    @{ Row = 82; Vals = @("Uruguay", 338, 0, 41, 296, 12, 0, 1) },
    @{ Row = 83; Vals = @("Republica de Macedonia", 329, 0, 12, 308, 4, 0, 9) },
    @{ Row = 84; Vals = @("Taiwan", 329, 7, 45, 279, 0, 0, 5) },
    @{ Row = 85; Vals = @("Kuwait", 317, 28, 80, 237, 13, 0, 0) },
    @{ Row = 92; Vals = @("Camerun", 233, 40, 10, 217, 0, 0, 6) },
    @{ Row = 97; Vals = @("Senegal", 190, 15, 45, 144, 0, 1, 1) },
    @{ Row = 98; Vals = @("Malta", 188, 19, 2, 186, 2, 0, 0) },
    @{ Row = 99; Vals = @("Cuba", 186, 0, 8, 172, 3, 0, 6) },
    @{ Row = 100; Vals = @("Costa de Marfil", 179, 0, 7, 171, 0, 0, 1) },
    @{ Row = 101; Vals = @("Uzbekistan", 173, 1, 8, 163, 8, 0, 2) },
    @{ Row = 102; Vals = @("Islas Feroe", 173, 4, 75, 98, 1, 0, 0) },
    @{ Row = 103; Vals = @("Honduras", 172, 31, 3, 159, 4, 3, 10) },
    @{ Row = 110; Vals = @("Brunei", 131, 2, 52, 78, 3, 0, 1) },
    @{ Row = 126; Vals = @("Kenia", 59, 0, 3, 55, 2, 0, 1) },
    @{ Row = 142; Vals = @("El Salvador", 32, 0, 0, 31, 1, 0, 1) },
    @{ Row = 143; Vals = @("Guam", 32, 0, 0, 31, 0, 0, 1) },
    @{ Row = 156; Vals = @("San Martin (Parte Holandesa)", 16, 10, 6, 9, 0, 1, 1) },
    @{ Row = 158; Vals = @("Bahamas", 15, 1, 1, 14, 0, 0, 0) },
    @{ Row = 160; Vals = @("Guinea Ecuatorial", 15, 0, 1, 14, 0, 0, 0) },
    @{ Row = 161; Vals = @("San Martin (Parte Francesa)", 15, 0, 2, 12, 0, 0, 1) },
    @{ Row = 162; Vals = @("Islas Caimanes", 14, 0, 0, 13, 0, 0, 1) },
    @{ Row = 163; Vals = @("Mongolia", 14, 2, 2, 12, 0, 0, 0) },
    @{ Row = 164; Vals = @("Santa Lucia", 13, 0, 1, 12, 0, 0, 0) },
    @{ Row = 165; Vals = @("Dominica", 12, 0, 0, 12, 0, 0, 0) },
    @{ Row = 166; Vals = @("Guyana", 12, 0, 0, 10, 0, 0, 2) },
    @{ Row = 167; Vals = @("Namibia", 11, 0, 2, 9, 0, 0, 0) },
    @{ Row = 168; Vals = @("Curazao", 11, 0, 2, 8, 0, 0, 1) },
    @{ Row = 171; Vals = @("Seychelles", 10, 0, 0, 10, 0, 0, 0) },
    @{ Row = 172; Vals = @("Libia", 10, 0, 1, 9, 0, 0, 0) },
    @{ Row = 174; Vals = @("Groenlandia", 10, 0, 2, 8, 0, 0, 0) },
    @{ Row = 176; Vals = @("Suazilandia", 9, 0, 0, 9, 0, 0, 0) },
    @{ Row = 177; Vals = @("Benin", 9, 0, 1, 8, 0, 0, 0) },
    @{ Row = 178; Vals = @("San Cristobal y Nieves", 8, 0, 0, 8, 0, 0, 0) },
    @{ Row = 180; Vals = @("Guinea-Bisau", 8, 0, 0, 8, 0, 0, 0) },
    @{ Row = 181; Vals = @("Zimbabue", 8, 0, 0, 7, 0, 0, 1) },
    @{ Row = 182; Vals = @("Antigua y Barbuda", 7, 0, 0, 7, 0, 0, 0) },
    @{ Row = 183; Vals = @("Republica del Chad", 7, 0, 0, 7, 0, 0, 0) },
    @{ Row = 185; Vals = @("Angola", 7, 0, 1, 4, 0, 0, 2) },
    @{ Row = 188; Vals = @("San Bartolome", 6, 0, 1, 5, 0, 0, 0) },
    @{ Row = 189; Vals = @("Cabo Verde", 6, 0, 0, 5, 0, 0, 1) },
    @{ Row = 191; Vals = @("Fiyi", 5, 0, 0, 5, 0, 0, 0) },
    @{ Row = 192; Vals = @("Montserrat", 5, 0, 0, 5, 0, 0, 0) },
    @{ Row = 193; Vals = @("Islas Turcas y Caicos", 5, 0, 0, 5, 0, 0, 0) },
    @{ Row = 198; Vals = @("Botsuana", 4, 0, 0, 3, 0, 0, 1) },
    @{ Row = 199; Vals = @("Gambia", 4, 0, 0, 3, 0, 0, 1) },
    @{ Row = 200; Vals = @("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 202; Vals = @("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 206; Vals = @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 207; Vals = @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 208; Vals = @("Sierra Leona", 1, 0, 0, 1, 0, 0, 0) }
)

foreach ($r in $rowsData) {
    $arr = New-Object 'object[,]' 1,8
    for ($i = 0; $i -lt 8; $i++) {
        $arr[0,$i] = $r.Vals[$i]
    }
    $rowNum = $r.Row
    $ws.Range("A${rowNum}:H${rowNum}").Value = $arr
}
